$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1249.75
$ws.Range("I29").Value = 999.7143
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 2999.1429
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = -2718.1429
$ws.Range("N29").Value = -9562
$ws.Range("H80").Value = 568.4761999999999
$ws.Range("I80").Value = 328.27274
$ws.Range("J80").Value = 832.7
$ws.Range("K80").Value = 984.81822
$ws.Range("L80").Value = 2498.1
$ws.Range("M80").Value = 13.18178
$ws.Range("N80").Value = -4494.1
$ws.Range("H83").Value = 568.4761999999999
$ws.Range("I83").Value = 328.27274
$ws.Range("J83").Value = 832.7
$ws.Range("K83").Value = 2954.45466
$ws.Range("L83").Value = 7494.3
$ws.Range("M83").Value = 2037.54534
$ws.Range("N83").Value = -17478.3
$ws.Range("H86").Value = 1917.3334
$ws.Range("I86").Value = 1901.25
$ws.Range("J86").Value = 1925.375
$ws.Range("K86").Value = 1901.25
$ws.Range("L86").Value = 1925.375
$ws.Range("M86").Value = -778.25
$ws.Range("N86").Value = -4171.375
$ws.Range("H88").Value = 14150
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 15533.333
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 15533.333
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -16345.333
$ws.Range("H89").Value = 1917.3334
$ws.Range("I89").Value = 1901.25
$ws.Range("J89").Value = 1925.375
$ws.Range("K89").Value = 9506.25
$ws.Range("L89").Value = 9626.875
$ws.Range("M89").Value = -3890.25
$ws.Range("N89").Value = -20858.875
$ws.Range("H91").Value = 14150
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 15533.333
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 15533.333
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -18341.333
$ws.Range("H97").Value = 1536.6666
$ws.Range("I97").Value = 1750
$ws.Range("J97").Value = 1110
$ws.Range("K97").Value = 5250
$ws.Range("L97").Value = 3330
$ws.Range("M97").Value = -4754
$ws.Range("N97").Value = -4322
$ws.Range("H138").Value = 2738.9597
$ws.Range("I138").Value = 990.8
$ws.Range("J138").Value = 3051.1309
$ws.Range("K138").Value = 2972.4
$ws.Range("L138").Value = 9153.3927
$ws.Range("M138").Value = 2167.6
$ws.Range("N138").Value = -19433.3927

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1640.1072
$ws.Range("I74").Value = 891.2222
$ws.Range("J74").Value = 2988.1
$ws.Range("K74").Value = 891.2222
$ws.Range("L74").Value = 2988.1
$ws.Range("M74").Value = -17.22220000000004
$ws.Range("N74").Value = -4736.1
$ws.Range("H77").Value = 1640.1072
$ws.Range("I77").Value = 891.2222
$ws.Range("J77").Value = 2988.1
$ws.Range("K77").Value = 4456.111
$ws.Range("L77").Value = 14940.5
$ws.Range("M77").Value = -88.11099999999988
$ws.Range("N77").Value = -23676.5
$ws.Range("H97").Value = 1238.28
$ws.Range("I97").Value = 804.0526
$ws.Range("K97").Value = 804.0526
$ws.Range("M97").Value = -308.0526
$ws.Range("H110").Value = 1008.36365
$ws.Range("I110").Value = 1051.7142
$ws.Range("J110").Value = 932.5
$ws.Range("K110").Value = 1051.7142
$ws.Range("L110").Value = 932.5
$ws.Range("M110").Value = 993.2858000000001
$ws.Range("N110").Value = -5022.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 28440
$ws.Range("J51").Value = 28440
$ws.Range("L51").Value = 28440
$ws.Range("N51").Value = -29422
$ws.Range("H86").Value = 2385.7856
$ws.Range("I86").Value = 1641.2
$ws.Range("J86").Value = 2799.4443
$ws.Range("K86").Value = 1641.2
$ws.Range("L86").Value = 2799.4443
$ws.Range("M86").Value = -518.2
$ws.Range("N86").Value = -5045.4443
$ws.Range("H89").Value = 2385.7856
$ws.Range("I89").Value = 1641.2
$ws.Range("J89").Value = 2799.4443
$ws.Range("K89").Value = 8206
$ws.Range("L89").Value = 13997.2215
$ws.Range("M89").Value = -2590
$ws.Range("N89").Value = -25229.2215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3032.861
$ws.Range("I31").Value = 1248.1
$ws.Range("J31").Value = 5263.8125
$ws.Range("K31").Value = 1248.1
$ws.Range("L31").Value = 5263.8125
$ws.Range("M31").Value = -953.0999999999999
$ws.Range("N31").Value = -5853.8125
$ws.Range("H34").Value = 3032.861
$ws.Range("I34").Value = 1248.1
$ws.Range("J34").Value = 5263.8125
$ws.Range("K34").Value = 1248.1
$ws.Range("L34").Value = 5263.8125
$ws.Range("M34").Value = -1046.1
$ws.Range("N34").Value = -5667.8125
$ws.Range("H122").Value = 4214.4
$ws.Range("I122").Value = 1518
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 4554
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2104
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1500
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1500
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 4500
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -4838
$ws.Range("H30").Value = 1500
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1500
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 4500
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -4704
$ws.Range("H49").Value = 2525.75
$ws.Range("I49").Value = 2525.75
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 7577.25
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -7421.25
$ws.Range("N49").ClearContents()
$ws.Range("H140").Value = 2597.524
$ws.Range("I140").Value = 2891.0588
$ws.Range("J140").Value = 1350
$ws.Range("K140").Value = 8673.1764
$ws.Range("L140").Value = 4050
$ws.Range("M140").Value = -3493.1764
$ws.Range("N140").Value = -14410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5507.9473
$ws.Range("I70").Value = 5125.35
$ws.Range("J70").Value = 6408.1763
$ws.Range("K70").Value = 5125.35
$ws.Range("L70").Value = 6408.1763
$ws.Range("M70").Value = -4855.35
$ws.Range("N70").Value = -6948.1763
$ws.Range("H73").Value = 5507.9473
$ws.Range("I73").Value = 5125.35
$ws.Range("J73").Value = 6408.1763
$ws.Range("K73").Value = 5125.35
$ws.Range("L73").Value = 6408.1763
$ws.Range("M73").Value = -4189.35
$ws.Range("N73").Value = -8280.176299999999
$ws.Range("H80").Value = 50004400
$ws.Range("I80").Value = 125002500
$ws.Range("J80").Value = 5666.6665
$ws.Range("K80").Value = 125002500
$ws.Range("L80").Value = 5666.6665
$ws.Range("M80").Value = -125001502
$ws.Range("N80").Value = -7662.6665
$ws.Range("H83").Value = 50004400
$ws.Range("I83").Value = 125002500
$ws.Range("J83").Value = 5666.6665
$ws.Range("K83").Value = 625012500
$ws.Range("L83").Value = 28333.3325
$ws.Range("M83").Value = -625007508
$ws.Range("N83").Value = -38317.3325
$ws.Range("H132").Value = 3894.5715
$ws.Range("I132").Value = 2761.75
$ws.Range("J132").Value = 7519.6
$ws.Range("K132").Value = 8285.25
$ws.Range("L132").Value = 22558.8
$ws.Range("M132").Value = -5755.25
$ws.Range("N132").Value = -27618.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5008.4814
$ws.Range("I82").Value = 6169.3887
$ws.Range("K82").Value = 6169.3887
$ws.Range("M82").Value = -5808.3887
$ws.Range("H85").Value = 5008.4814
$ws.Range("I85").Value = 6169.3887
$ws.Range("K85").Value = 6169.3887
$ws.Range("M85").Value = -4921.3887
$ws.Range("H93").Value = 8549615
$ws.Range("I93").Value = 13890749
$ws.Range("J93").Value = 3800
$ws.Range("K93").Value = 13890749
$ws.Range("L93").Value = 3800
$ws.Range("M93").Value = -13889501
$ws.Range("N93").Value = -6296

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2702836
$ws.Range("I96").Value = 167999.67
$ws.Range("J96").Value = 3970254
$ws.Range("K96").Value = 167999.67
$ws.Range("L96").Value = 3970254
$ws.Range("M96").Value = -166626.67
